$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 11.83

$ws.Range("D3").Value = 10
$ws.Range("E3").Value = 10.17

$ws.Range("C4").Value = 10
$ws.Range("E4").Value = 9.76
$ws.Range("F4").Value = 10.31

$ws.Range("B5").Value = 8.17
$ws.Range("C5").Value = 9.83
$ws.Range("D5").Value = 10.24
$ws.Range("F5").Value = 11.18
$ws.Range("G5").Value = 10.89

$ws.Range("D6").Value = 9.69
$ws.Range("E6").Value = 8.82
$ws.Range("G6").Value = 11.21
$ws.Range("H6").Value = 12.75
$ws.Range("J6").Value = 3.33

$ws.Range("E7").Value = 9.11
$ws.Range("F7").Value = 8.79
$ws.Range("H7").Value = 6.67

$ws.Range("F8").Value = 7.25
$ws.Range("G8").Value = 13.33

$ws.Range("F10").Value = 16.67
